$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9931485056877136
$ws.Range("B1").Value = 1.444755554199219
$ws.Range("C1").Value = 2.963318109512329
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.508775949478149
